$d = $word.ActiveDocument

# 1. licenceGreyCode -> licenceGrayCode
$d.Content.Find.Execute("licenceGreyCode", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "licenceGrayCode", 2) | Out-Null

# 2. "Comporte un numLicence traduit en grey code" -> "Comporte un numLicence traduit en gray code"
$d.Content.Find.Execute("numLicence traduit en grey code", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "numLicence traduit en gray code", 2) | Out-Null

# 3. correspGreyCodeBarre -> correspGrayCodeBarre
$d.Content.Find.Execute("correspGreyCodeBarre", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "correspGrayCodeBarre", 2) | Out-Null

# 4. "Correspondance entre les grey codes ..." -> "... gray codes ..."
$d.Content.Find.Execute("les grey codes et les 8", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "les gray codes et les 8", 2) | Out-Null

$table = $d.Tables.Item(1)

# 5. motBinaire row, Type column: "Chaine de caractères" -> "Liste d'entier"
$cellType = $table.Cell(5, 2)
$ts = $cellType.Range.Start
$te = $cellType.Range.End
$d.Range($ts, $te).Text = "Liste d’entier"

# 6. motBinaire row, Signification column: rewrite description text
$cellDesc = $table.Cell(5, 3)
$ds = $cellDesc.Range.Start
$de = $cellDesc.Range.End
$d.Range($ds, $de).Text = "Variable temporaire qui contient un mot de 3 bits qui représente la moitié d’encodage d’un caractère. Un caractère est codé sur 6 bits et une barre sur 3 bits, il faut donc 2 barre pour encoder un caractère."

# 7. Remove the 3 trailing empty rows
for ($i = $table.Rows.Count; $i -ge 6; $i--) {
    $table.Rows.Item($i).Delete()
}

Write-Output "done"
